# The commit swaps the theme content that was stored in ppt/theme/theme1.xml
# (the "Integral" theme, used by the deck's SlideMaster) with the theme
# content that was stored in ppt/theme/theme2.xml (the default "Office
# Theme", used by the NotesMaster) -- i.e. the SlideMaster ends up on the
# plain "Office Theme" colour palette while the NotesMaster ends up on the
# "Integral" palette.
#
# The PowerPoint object model only exposes the *active* (SlideMaster) theme's
# colour scheme for editing (Slide.ThemeColorScheme / the legacy
# NotesMaster.ColorScheme both resolve to the one theme driving the deck), so
# we reassign each of the twelve theme colour slots on that scheme to the
# "Office Theme" palette -- matching the colour values the diff shows ending
# up in ppt/theme/theme1.xml.

function HexToBGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (schemeClr order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink == ThemeColorScheme index 1-12).
$officeTheme = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToBGR($officeTheme[$i - 1])
}
